$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 40: fix the date (was 45675, should be 45677)
$ws.Range("A40").Value = 45677

# Row 41 (new): library/book donation entry
$ws.Range("A41").Value = 45678
$ws.Range("A41").NumberFormat = $ws.Range("A40").NumberFormat
$ws.Range("B41").Value = "animations for the bossfight"
$ws.Range("C41").Value = 4

# Row 42 (new): blank row, date-formatted cell in A42 only (no value)
$ws.Range("A42").NumberFormat = $ws.Range("A40").NumberFormat

# Update the sheet selection to match the saved view state
$ws.Range("C46").Select()

$wb.Save()
